$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo in the report title (A1): "Roport" -> "Raport"
$ws.Range("A1").Value = "Raport zaangażowania"

# 2. Add the new engagement data (08.05.2025) in columns B:D, rows 9-17
$ws.Range("B9").Value = 45785
$ws.Range("C9").Value = "CustomDummyTest.cs"
$ws.Range("D9").Value = 37

$ws.Range("C10").Value = "CustomFakeTest.cs"
$ws.Range("D10").Value = 41

$ws.Range("C11").Value = "CustomMockTest.cs"
$ws.Range("D11").Value = 43

$ws.Range("C12").Value = "CustomSpyTest.cs"
$ws.Range("D12").Value = 44

$ws.Range("C13").Value = "CustomStubTest.cs"
$ws.Range("D13").Value = 36

$ws.Range("C14").Value = "MoqDummyTest.cs"
$ws.Range("D14").Value = 22

$ws.Range("C15").Value = "MoqMockTest.cs"
$ws.Range("D15").Value = 23

$ws.Range("C16").Value = "MoqStubTest.cs"
$ws.Range("D16").Value = 26

$ws.Range("C17").Value = "TestyJednostkoweBLL.csproj"
$ws.Range("D17").Value = 30

# 3. Match the existing date-column formatting (centered short date, like column K)
$ws.Range("K8").Copy()
$ws.Range("B9:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Merge the new date cell down the block, like the other date columns
$ws.Range("B9:B17").Merge()

# 5. Widen column C so the new, longer file names fit
$ws.Columns("C").ColumnWidth = 25.85546875

# 6. Move the active selection to reflect where the new data was entered
$ws.Range("E25").Select()
